$d = $word.ActiveDocument

function Find-Paragraph($doc, $pattern) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like $pattern) {
            return $p
        }
    }
    return $null
}

function Find-ParagraphIndex($doc, $pattern) {
    $idx = 0
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text -like $pattern) {
            return $idx
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1) Add a new "Driver Controlled Enhancements" section just before
#    "Engineering Portfolio References": a bold heading paragraph plus
#    a body paragraph describing the new PIDF-loop / incremental-move
#    driver controls.
# ------------------------------------------------------------------

$keyAlgBodyIdx = Find-ParagraphIndex $d "*tapes were too small.*"
$keyAlgBody = $d.Paragraphs($keyAlgBodyIdx)
$keyAlgBody.Range.InsertParagraphAfter()

$headingBlank = $d.Paragraphs($keyAlgBodyIdx + 1)
$headingXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Driver Controlled Enhancements</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$headingBlank.Range.InsertXML($headingXml)

$headingIdx = Find-ParagraphIndex $d "Driver Controlled Enhancements*"
$headingPara = $d.Paragraphs($headingIdx)
$headingPara.Range.InsertParagraphAfter()

$pidfBlank = $d.Paragraphs($headingIdx + 1)
$pidfXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t xml:space="preserve">To move the arm up at a specific angle, we are using a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>pidf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> loop. When we press the “x” button on gamepad, the arm goes to the target position, making it easier to keep the arm up. In the chassis driver controls, we also added moving buttons that move the robot in small increments. This helps us be more accurate when in front of the backdrop and going through the truss.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pidfBlank.Range.InsertXML($pidfXml)

# ------------------------------------------------------------------
# 2) Add a new "Autonomous Diagram" section right after the
#    "...sections 4-6 and 10-12." paragraph (and before the existing
#    trailing empty paragraph): a bold heading paragraph plus an empty
#    bold paragraph.
# ------------------------------------------------------------------

$refsIdx = Find-ParagraphIndex $d "*sections 4-6 and 10-12.*"
$refsPara = $d.Paragraphs($refsIdx)
$refsPara.Range.InsertParagraphAfter()

$diagramBlank = $d.Paragraphs($refsIdx + 1)
$diagramXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Autonomous Diagram</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$diagramBlank.Range.InsertXML($diagramXml)

$diagramIdx = Find-ParagraphIndex $d "Autonomous Diagram*"
$diagramPara = $d.Paragraphs($diagramIdx)
$diagramPara.Range.InsertParagraphAfter()

$emptyBoldBlank = $d.Paragraphs($diagramIdx + 1)
$emptyBoldXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$emptyBoldBlank.Range.InsertXML($emptyBoldXml)
